$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 267; existing rows 267..278 shift down to 268..279.
$ws.Rows.Item(267).Insert()

# Populate the newly inserted row 267 with the new weekly record.
$ws.Cells.Item(267, 1).Value = 7
$ws.Cells.Item(267, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(267, 3).Value = "Ñuble"
$ws.Cells.Item(267, 4).Value = 44509
$ws.Cells.Item(267, 5).Value = 16
$ws.Cells.Item(267, 6).Value = 100114014
$ws.Cells.Item(267, 7).Value = "Betarraga"
$ws.Cells.Item(267, 8).Value = "Sin especificar"
$ws.Cells.Item(267, 9).Value = "Primera"
$ws.Cells.Item(267, 10).Value = 240
$ws.Cells.Item(267, 11).Value = 700
$ws.Cells.Item(267, 12).Value = 800
$ws.Cells.Item(267, 13).Value = 750
$ws.Cells.Item(267, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(267, 15).Value = "Región del Maule"
$ws.Cells.Item(267, 16).Value = 150
$ws.Cells.Item(267, 17).Value = 5
$ws.Cells.Item(267, 18).Value = "Hortaliza"
